$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 369; this shifts the existing rows
# 369:489 down to 370:490, preserving all their data/formatting.
$ws.Rows(369).Insert()

# Populate the newly inserted row 369 with the new data record.
$ws.Cells.Item(369, 1).Value = 3
$ws.Cells.Item(369, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(369, 3).Value = "Coquimbo"
$ws.Cells.Item(369, 4).Value = 44588
$ws.Cells.Item(369, 5).Value = 5
$ws.Cells.Item(369, 6).Value = 100112006
$ws.Cells.Item(369, 7).Value = "Repollo"
$ws.Cells.Item(369, 8).Value = "Crespo record"
$ws.Cells.Item(369, 9).Value = "Primera"
$ws.Cells.Item(369, 10).Value = 1550
$ws.Cells.Item(369, 11).Value = 1000
$ws.Cells.Item(369, 12).Value = 1100
$ws.Cells.Item(369, 13).Value = 1052
$ws.Cells.Item(369, 14).Value = "$/unidad"
$ws.Cells.Item(369, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(369, 16).Value = 1052
$ws.Cells.Item(369, 17).Value = 1
$ws.Cells.Item(369, 18).Value = "Hortaliza"
